$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "'36.885.23"
$ws.Range("E2").Formula = "'  -0.41%  "

# Row 3
$ws.Range("D3").Formula = "'2.051.66"
$ws.Range("E3").Formula = "'  +0.47%  "

# Row 4
$ws.Range("E4").Formula = "'  -0.23%  "

# Row 5
$ws.Range("D5").Formula = "'245.47"
$ws.Range("E5").Formula = "'  -1.10%  "

# Row 6
$ws.Range("D6").Formula = "'0.654"
$ws.Range("E6").Formula = "'  -1.49%  "

# Row 7
$ws.Range("B7").Formula = "'Solana"
$ws.Range("C7").Formula = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Formula = "'57.85"
$ws.Range("E7").Formula = "'  -1.91%  "

# Row 8
$ws.Range("B8").Formula = "'USDC"
$ws.Range("C8").Formula = "'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Formula = "'1.00"
$ws.Range("E8").Formula = "'  +0.01%  "

# Row 9
$ws.Range("D9").Formula = "'58.59"
$ws.Range("E9").Formula = "'  -2.60%  "

# Row 10
$ws.Range("D10").Formula = "'0.368"
$ws.Range("E10").Formula = "'  -3.91%  "

# Row 11
$ws.Range("D11").Formula = "'0.0776"
$ws.Range("E11").Formula = "'  -1.22%  "

# Row 12
$ws.Range("E12").Formula = "'  +1.84%  "

# Row 13
$ws.Range("D13").Formula = "'15.21"
$ws.Range("E13").Formula = "'  -3.61%  "

# Row 14
$ws.Range("D14").Formula = "'0.872"
$ws.Range("E14").Formula = "'  +5.22%  "

# Row 15
$ws.Range("D15").Formula = "'2.348.08"
$ws.Range("E15").Formula = "'  +0.41%  "

# Row 16
$ws.Range("D16").Formula = "'5.60"
$ws.Range("E16").Formula = "'  -2.45%  "

# Row 17
$ws.Range("D17").Formula = "'2.062.38"
$ws.Range("E17").Formula = "'  +1.19%  "

# Row 18
$ws.Range("B18").Formula = "'WrappedBTC"
$ws.Range("C18").Formula = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Formula = "'36.805.17"
$ws.Range("E18").Formula = "'  -0.66%  "

# Row 19
$ws.Range("B19").Formula = "'Avalanche"
$ws.Range("C19").Formula = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Formula = "'17.59"
$ws.Range("E19").Formula = "'  -5.15%  "

# Row 20
$ws.Range("D20").Formula = "'73.27"
$ws.Range("E20").Formula = "'  -1.99%  "

# Row 21
$ws.Range("D21").Formula = "'0.0₃0885"
$ws.Range("E21").Formula = "'  -1.31%  "

# Row 22
$ws.Range("D22").Formula = "'5.38"
$ws.Range("E22").Formula = "'  +0.84%  "

# Row 23
$ws.Range("D23").Formula = "'236.12"
$ws.Range("E23").Formula = "'  -0.21%  "

# Row 24
$ws.Range("E24").Formula = "'  +0.00%  "

# Row 25
$ws.Range("D25").Formula = "'2.45"
$ws.Range("E25").Formula = "'  +1.94%  "

# Row 26
$ws.Range("D26").Formula = "'10.40"
$ws.Range("E26").Formula = "'  +11.23%  "

# Row 27
$ws.Range("D27").Formula = "'2.23"
$ws.Range("E27").Formula = "'  +2.33%  "

# Row 28
$ws.Range("D28").Formula = "'168.65"
$ws.Range("E28").Formula = "'  -0.19%  "

# Row 29
$ws.Range("D29").Formula = "'19.96"
$ws.Range("E29").Formula = "'  -0.37%  "

# Row 30
$ws.Range("D30").Formula = "'5.52"
$ws.Range("E30").Formula = "'  +15.88%  "

# Row 31
$ws.Range("D31").Formula = "'0.124"
$ws.Range("E31").Formula = "'  -1.07%  "

# Row 32
$ws.Range("D32").Formula = "'1.11"
$ws.Range("E32").Formula = "'  -1.95%  "

# Row 33
$ws.Range("D33").Formula = "'4.86"
$ws.Range("E33").Formula = "'  +7.44%  "

# Row 34
$ws.Range("D34").Formula = "'0.0612"
$ws.Range("E34").Formula = "'  -2.11%  "

# Row 35
$ws.Range("B35").Formula = "'LidoDAOToken"
$ws.Range("C35").Formula = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Formula = "'2.33"
$ws.Range("E35").Formula = "'  +5.40%  "

# Row 36
$ws.Range("B36").Formula = "'BinanceUSD"
$ws.Range("C36").Formula = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").Formula = "'1.00"
$ws.Range("E36").Formula = "'  -0.07%  "

# Row 37
$ws.Range("E37").Formula = "'  +4.36%  "

# Row 38
$ws.Range("E38").Formula = "'  -10.20%  "

# Row 39
$ws.Range("E39").Formula = "'  -2.06%  "

# Row 40
$ws.Range("B40").Formula = "'THORChain"
$ws.Range("C40").Formula = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").Formula = "'4.98"
$ws.Range("E40").Formula = "'  -4.25%  "

# Row 41
$ws.Range("B41").Formula = "'HuobiToken"
$ws.Range("C41").Formula = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Formula = "'2.99"
$ws.Range("E41").Formula = "'  -6.01%  "

# Row 42
$ws.Range("D42").Formula = "'0.0222"
$ws.Range("E42").Formula = "'  +0.44%  "

# Row 43
$ws.Range("D43").Formula = "'1.15"
$ws.Range("E43").Formula = "'  +1.48%  "

# Row 44
$ws.Range("B44").Formula = "'InjectiveProtocol"
$ws.Range("C44").Formula = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Formula = "'16.93"
$ws.Range("E44").Formula = "'  -2.60%  "

# Row 45
$ws.Range("B45").Formula = "'Cronos"
$ws.Range("C45").Formula = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Formula = "'0.0940"
$ws.Range("E45").Formula = "'  -11.89%  "

# Row 46
$ws.Range("B46").Formula = "'Aave"
$ws.Range("C46").Formula = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Formula = "'96.30"
$ws.Range("E46").Formula = "'  +0.46%  "

# Row 47
$ws.Range("D47").Formula = "'1.312.47"
$ws.Range("E47").Formula = "'  +1.97%  "

# Row 48
$ws.Range("D48").Formula = "'2.36"
$ws.Range("E48").Formula = "'  -4.49%  "

# Row 49
$ws.Range("E49").Formula = "'  -2.05%  "

# Row 50
$ws.Range("D50").Formula = "'6.72"
$ws.Range("E50").Formula = "'  -1.05%  "

# Row 51
$ws.Range("D51").Formula = "'2.235.15"
$ws.Range("E51").Formula = "'  +0.59%  "
